# Apply the Jogos_da_Semana_FlashScore_2024-11-17.xlsx update:
#  - update a handful of odds values in row 3
#  - remove rows 6 and 7 (which shrinks the used range to A1:BD5)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update odds values on row 3
$ws.Range("I3").Value = 5.75
$ws.Range("J3").Value = 2.6
$ws.Range("X3").Value = 6.5
$ws.Range("AN3").Value = 3.4
$ws.Range("AV3").Value = 7

# Delete rows 6 and 7 entirely (row 7 first so row numbers don't shift
# before row 6 is removed)
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()
